# Update countries & provincias Spain
# Applies:
#  1) Re-sorted country labels (3 extra pairs swap position in the
#     "País" column because their underlying case totals crossed over)
#  2) Refreshed numeric statistics (Casos totales / Nuevos casos /
#     Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes)
#  3) Updated "datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-CountryLabel($row1, $row2) {
    $label1 = $ws.Cells.Item($row1, 1).Value2
    $label2 = $ws.Cells.Item($row2, 1).Value2
    $ws.Cells.Item($row1, 1).Value = $label2
    $ws.Cells.Item($row2, 1).Value = $label1
}

# --- 1) Swap labels where the underlying ranking order changed ---
Swap-CountryLabel 35 36    # Kuwait <-> Singapur
Swap-CountryLabel 202 203  # Fiyi <-> Dominica
Swap-CountryLabel 208 209  # Groenlandia <-> Islas Malvinas
Swap-CountryLabel 211 212  # Seychelles <-> Montserrat

# --- 2) Refresh numeric data (B=Casos totales, C=Nuevos casos,
#         D=Casos activos, E=Recuperados, F=Casos criticos,
#         G=Muertes hoy, H=Muertes) ---

# Row 6 (Rusia)
$ws.Cells.Item(6,2).Value = 620794
$ws.Cells.Item(6,3).Value = 6800
$ws.Cells.Item(6,4).Value = 384152
$ws.Cells.Item(6,5).Value = 227861
$ws.Cells.Item(6,7).Value = 176
$ws.Cells.Item(6,8).Value = 8781

# Row 7 (India)
$ws.Cells.Item(7,2).Value = 491741
$ws.Cells.Item(7,3).Value = 571
$ws.Cells.Item(7,4).Value = 285983
$ws.Cells.Item(7,5).Value = 190439
$ws.Cells.Item(7,7).Value = 11
$ws.Cells.Item(7,8).Value = 15319

# Row 20 (Banglades)
$ws.Cells.Item(20,2).Value = 130474
$ws.Cells.Item(20,3).Value = 3868
$ws.Cells.Item(20,4).Value = 53133
$ws.Cells.Item(20,5).Value = 75680
$ws.Cells.Item(20,7).Value = 40
$ws.Cells.Item(20,8).Value = 1661

# Row 35 (now Singapur after swap)
$ws.Cells.Item(35,2).Value = 42955
$ws.Cells.Item(35,3).Value = 219
$ws.Cells.Item(35,4).Value = 36604
$ws.Cells.Item(35,5).Value = 6325
$ws.Cells.Item(35,8).Value = 26

# Row 36 (now Kuwait after swap)
$ws.Cells.Item(36,2).Value = 42788
$ws.Cells.Item(36,4).Value = 33367
$ws.Cells.Item(36,5).Value = 9082
$ws.Cells.Item(36,8).Value = 339

# Row 41 (Polonia)
$ws.Cells.Item(41,2).Value = 33395
$ws.Cells.Item(41,3).Value = 276
$ws.Cells.Item(41,4).Value = 19218
$ws.Cells.Item(41,5).Value = 12748
$ws.Cells.Item(41,7).Value = 17
$ws.Cells.Item(41,8).Value = 1429

# Row 44 (Afganistan)
$ws.Cells.Item(44,2).Value = 30451
$ws.Cells.Item(44,3).Value = 276
$ws.Cells.Item(44,4).Value = 10306
$ws.Cells.Item(44,5).Value = 19462
$ws.Cells.Item(44,7).Value = 8
$ws.Cells.Item(44,8).Value = 683

# Row 110 (Estonia)
$ws.Cells.Item(110,2).Value = 1986
$ws.Cells.Item(110,3).Value = 2
$ws.Cells.Item(110,4).Value = 1797
$ws.Cells.Item(110,5).Value = 120

# Row 114 (Lituania)
$ws.Cells.Item(114,2).Value = 1808
$ws.Cells.Item(114,3).Value = 2
$ws.Cells.Item(114,4).Value = 1501
$ws.Cells.Item(114,5).Value = 229

# Row 116 (Eslovaquia)
$ws.Cells.Item(116,2).Value = 1643
$ws.Cells.Item(116,3).Value = 13
$ws.Cells.Item(116,4).Value = 1455
$ws.Cells.Item(116,5).Value = 160

# Row 156 (Montenegro)
$ws.Cells.Item(156,2).Value = 424
$ws.Cells.Item(156,3).Value = 10
$ws.Cells.Item(156,5).Value = 100

# Row 206 (Islas Turcas y Caicos)
$ws.Cells.Item(206,2).Value = 16
$ws.Cells.Item(206,3).Value = 1
$ws.Cells.Item(206,5).Value = 4

# Row 211 (now Montserrat after swap)
$ws.Cells.Item(211,4).Value = 10
$ws.Cells.Item(211,8).Value = 1

# Row 212 (now Seychelles after swap)
$ws.Cells.Item(212,4).Value = 11
$ws.Cells.Item(212,8).Value = 0

# --- 3) Update the "datos actualizados" timestamp ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 26 de Junio de 2020 a las 10:43"
